# Reposition the four pictures on slide 17 (the "Change in Cognitive Tests
# Performance of Thresholds" slide) slightly lower, matching the layout
# adjustment made when the PDF export of the deck was added.
#
# Only the vertical offset (Top) changes for each picture; the horizontal
# offset (Left) and size stay the same.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# Picture 2 (id 2050) - left-hand chart image
$s.Shapes.Item(3).Top = 107.16833021653544

# Picture 3 (id 4) - right-hand chart image
$s.Shapes.Item(4).Top = 107.93336958661418

# Picture 11 (id 12) - small cropped legend/callout over the right chart
$s.Shapes.Item(5).Top = 135.44574803149607

# Picture 2 (id 13) - small cropped legend/callout over the left chart
$s.Shapes.Item(6).Top = 124.95147982283464
